$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRICES")

# Copy the date formatting of the last existing date cell (A20) down onto
# the two new date cells (A21, A22) so they share the same style index
# instead of minting a new number-format style.
$ws.Range("A20").Copy()
$ws.Range("A21:A22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 21 - 2025-03-29
$ws.Cells.Item(21, 1).Value = 45745
$ws.Cells.Item(21, 2).Value = 23.88
$ws.Cells.Item(21, 3).Value = 4.5
$ws.Cells.Item(21, 4).Value = 2.5
$ws.Cells.Item(21, 5).Value = 11
$ws.Cells.Item(21, 6).Value = 4.2
$ws.Cells.Item(21, 7).Value = 5.7
$ws.Cells.Item(21, 8).Value = 1.8
$ws.Cells.Item(21, 9).Value = 1.795
$ws.Cells.Item(21, 10).Value = 89.95
$ws.Cells.Item(21, 11).Value = 24.95
$ws.Cells.Item(21, 12).Value = 589.95
$ws.Cells.Item(21, 13).Value = 49.95
$ws.Cells.Item(21, 14).Value = 134.95
$ws.Cells.Item(21, 15).Value = 71.95
$ws.Cells.Item(21, 16).Value = 46.72

# Row 22 - 2025-04-12
$ws.Cells.Item(22, 1).Value = 45759
$ws.Cells.Item(22, 2).Value = 23.92
$ws.Cells.Item(22, 3).Value = 4.5
$ws.Cells.Item(22, 4).Value = 3
$ws.Cells.Item(22, 5).Value = 11
$ws.Cells.Item(22, 6).Value = 4.2
$ws.Cells.Item(22, 7).Value = 6.2
$ws.Cells.Item(22, 8).Value = 1.8
$ws.Cells.Item(22, 9).Value = 1.699
$ws.Cells.Item(22, 10).Value = 95.95
$ws.Cells.Item(22, 11).Value = 16.95
$ws.Cells.Item(22, 12).Value = 589.95
$ws.Cells.Item(22, 13).Value = 49.95
$ws.Cells.Item(22, 14).Value = 154.95
$ws.Cells.Item(22, 15).Value = 74.95
$ws.Cells.Item(22, 16).Value = 44.18

# Match the saved cursor position recorded in the diff.
$ws.Range("P26").Select()
